$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.527.14"
$ws.Range("E2").Value = "  +0.34%  "

# Row 3
$ws.Range("D3").Value = "1.839.49"
$ws.Range("E3").Value = "  -0.06%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.00%  "

# Row 6
$ws.Range("E6").Value = "  +0.01%  "

# Row 7
$ws.Range("E7").Value = "  +0.37%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3192"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.58%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06795"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.20%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.78"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.54%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7844"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.27%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07757"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.89%  "

# Row 13
$ws.Range("D13").Value = "1.832.50"
$ws.Range("E13").Value = "  -0.85%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.74%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.015"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.26%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.02%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.70%  "

# Row 18
$ws.Range("E18").Value = "  +0.00%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007946"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.21%  "

# Row 20
$ws.Range("D20").Value = "26.540.92"
$ws.Range("E20").Value = "  +0.24%  "

# Row 21
$ws.Range("D21").Value = "2.069.11"
$ws.Range("E21").Value = "  -0.28%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.629"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.22%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.975"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.51%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.341"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.34%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.89%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.198"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.47%  "

# Row 27
$ws.Range("E27").Value = "  +1.30%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.04%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.37%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.169"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.43%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08705"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.62%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.082"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.57%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04878"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.29%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7297"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.66%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.137"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.08%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.860"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.51%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.094"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.96%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.246"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.35%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01756"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.68%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4816"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.45%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8947"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.56%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.75%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.938"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.68%  "

# Row 44
$ws.Range("E44").Value = "  +0.06%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.663"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.35%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4178"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.95%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.49%  "

# Row 48
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05853"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.06%  "

# Row 49
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1234"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.75%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.16%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8922"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.04%  "
